$d = $word.ActiveDocument

# Step 1: delete the "License Information" (Heading2) paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "License Information*") {
        $p.Range.Delete()
        break
    }
}

# Step 2: delete the "This PDF version is provided under the same license." paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "This PDF version is provided under the same license.*") {
        $p.Range.Delete()
        break
    }
}

# Step 3: locate the big license paragraph (the one that still starts with the
# bold "Questions de Traduction (unfoldingWord)" run followed by " (French) is based on").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Questions de Traduction (unfoldingWord) (French)*") {
        $target = $p
        break
    }
}

# Step 3a: change the bold run's text. We insert the new text right before the
# end of the old matched text and then delete the old text (rather than doing a
# straight Find/Replace or Range.Text= at the very start of the paragraph) so
# that the pre-existing leading empty run is not swallowed by the edit.
$rng = $target.Range
$rng.Find.Execute("Questions de Traduction (unfoldingWord)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oldStart = $rng.Start
$oldEnd = $rng.End
$insAfter = $d.Range($oldEnd, $oldEnd)
$insAfter.InsertBefore("unfoldingWord® Translation Questions")
$oldRange = $d.Range($oldStart, $oldEnd)
$oldRange.Delete()

# Step 3b: replace everything after the (new) bold run up to, but excluding,
# the trailing empty run before the paragraph mark, with the new text. Again
# we insert-then-delete (instead of a direct Range.Text= / Find replace) so
# that formatting is picked up from the non-bold tail (the old trailing
# period) rather than from the bold run that now precedes the insertion
# point, and so the single trailing empty run survives intact.
$rng2 = $target.Range
$rng2.Find.Execute("unfoldingWord® Translation Questions", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldEnd = $rng2.End

$oldRestStart = $boldEnd
$oldRestEnd = $target.Range.End - 1

$newText = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. " + `
    "unfoldingWord® Translation Questions" + `
    " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from " + `
    "unfoldingWord® Translation Questions" + `
    " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$insPoint = $d.Range($oldRestEnd, $oldRestEnd)
$insPoint.InsertBefore($newText)

$oldRestRange = $d.Range($oldRestStart, $oldRestEnd)
$oldRestRange.Delete()

Write-Output "done"
